# Add the new "2022-Q4" quarterly sheet (placed right after "总计") and
# refresh the "总计" summary sheet with the new quarter's row, pushing all
# the pre-existing quarters' rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after the "总计" sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Header row (row 1), columns B..H - column A header is blank, like the
# other quarter sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$q4Sheet.Range("B1:H1").Font.Bold = $true

# Force columns B (code) and D:G (amount/position/ratio/value, which look
# numeric but must stay text like the source data) to text before writing,
# so Excel does not silently coerce "011093" -> 11093 etc.
$q4Sheet.Range("B2:B10").NumberFormat = "@"
$q4Sheet.Range("D2:G10").NumberFormat = "@"

# Fund rows for 2022-Q4.
$rows = @(
    @("011093", "永赢宏泽一年定期开放灵活配置混合", "14.98", "48.20", "0.57", "0.0854", 6),
    @("160135", "南方中证高铁产业指数（LOF）",       "1.78",  "94.92", "2.47", "0.0440", 9),
    @("002210", "创金合信量化多因子股票A",           "3.02",  "93.22", "1.19", "0.0359", 8),
    @("006836", "永赢惠泽一年定期开放灵活配置混合", "3.73",  "48.14", "0.57", "0.0213", 6),
    @("003865", "创金合信量化多因子股票C",           "1.68",  "93.22", "1.19", "0.0200", 8),
    @("160639", "鹏华中证高铁产业指数（LOF）A",       "0.75",  "94.60", "2.45", "0.0184", 9),
    @("000892", "九泰天宝灵活配置混合A",             "0.06",  "94.55", "4.68", "0.0028", 6),
    @("015678", "鹏华中证高铁产业指数（LOF）C",       "0.06",  "94.60", "2.45", "0.0015", 9),
    @("002028", "九泰天宝灵活配置混合C",             "0.00",  "94.55", "4.68", 0,        6)
)

$r = 2
foreach ($row in $rows) {
    $q4Sheet.Cells.Item($r, 1).Value = $r - 2
    $q4Sheet.Cells.Item($r, 2).Value = $row[0]
    $q4Sheet.Cells.Item($r, 3).Value = $row[1]
    $q4Sheet.Cells.Item($r, 4).Value = $row[2]
    $q4Sheet.Cells.Item($r, 5).Value = $row[3]
    $q4Sheet.Cells.Item($r, 6).Value = $row[4]
    $q4Sheet.Cells.Item($r, 7).Value = $row[5]
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q4 and
#    push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Copy formatting of column A from the row below (the style that marks the
# numeric index column) onto the freshly inserted row, then clear the
# inherited formatting on B:D so the new row matches the plain look of the
# other data rows.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 0.23
